# Update Turkey macro dataset (module1: product, module2: macro) with the
# latest released figures. For each indicator row: "Last" (C) becomes the
# newly released reading, "Previous" (D) becomes the prior "Last" reading,
# and the announcement date (H) advances to the new release month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overview - Currency
$ws.Range("C2").Value = "34.35"
$ws.Range("D2").Value = "34.35"

# Overview - Stock Market
$ws.Range("C3").Value = "9300"
$ws.Range("D3").Value = "9227"

# Overview - Unemployment Rate
$ws.Range("C6").Value = "8.6"
$ws.Range("D6").Value = "8.6"
$ws.Range("H6").Value = "Sep/24"

# Overview - Current Account
$ws.Range("C11").Value = "2988"
$ws.Range("D11").Value = "4324"
$ws.Range("H11").Value = "Sep/24"

# Overview - Retail Sales MoM
$ws.Range("C18").Value = "2.3"
$ws.Range("D18").Value = "2.3"
$ws.Range("H18").Value = "Sep/24"

# Labour - Unemployment Rate
$ws.Range("C34").Value = "8.6"
$ws.Range("D34").Value = "8.6"
$ws.Range("H34").Value = "Sep/24"

# Labour - Employed Persons
$ws.Range("C35").Value = "32823"
$ws.Range("D35").Value = "32728"
$ws.Range("H35").Value = "Sep/24"

# Labour - Unemployed Persons
$ws.Range("C36").Value = "3100"
$ws.Range("D36").Value = "3065"
$ws.Range("H36").Value = "Sep/24"

# Labour - Labor Force Participation Rate
$ws.Range("C37").Value = "54.4"
$ws.Range("D37").Value = "54.2"
$ws.Range("H37").Value = "Sep/24"

# Labour - Youth Unemployment Rate
$ws.Range("C38").Value = "17.1"
$ws.Range("D38").Value = "16.8"
$ws.Range("H38").Value = "Sep/24"

# Labour - Employment Rate (date only)
$ws.Range("H46").Value = "Sep/24"

# Trade - Current Account
$ws.Range("C78").Value = "2988"
$ws.Range("D78").Value = "4324"
$ws.Range("H78").Value = "Sep/24"

# Trade - Capital Flows
$ws.Range("C84").Value = "-3263"
$ws.Range("D84").Value = "630"
$ws.Range("H84").Value = "Sep/24"

# Trade - Foreign Direct Investment
$ws.Range("C85").Value = "441"
$ws.Range("D85").Value = "304"
$ws.Range("H85").Value = "Sep/24"

# Trade - Remittances
$ws.Range("C86").Value = "14"
$ws.Range("D86").Value = "14"
$ws.Range("H86").Value = "Sep/24"

# Trade - Crude Oil Production
$ws.Range("C91").Value = "103"
$ws.Range("D91").Value = "102"
$ws.Range("H91").Value = "Jul/24"

# Business - Industrial Production
$ws.Range("C108").Value = "-2.4"
$ws.Range("D108").Value = "-5.2"
$ws.Range("H108").Value = "Sep/24"

# Business - Industrial Production MoM
$ws.Range("C109").Value = "1.6"
$ws.Range("D109").Value = "-1.6"
$ws.Range("H109").Value = "Sep/24"

# Business - Manufacturing Production
$ws.Range("C110").Value = "-2.5"
$ws.Range("D110").Value = "-5.4"
$ws.Range("H110").Value = "Sep/24"

# Business - Car Production
$ws.Range("C114").Value = "121970"
$ws.Range("D114").Value = "123445"
$ws.Range("H114").Value = "Oct/24"

# Business - Total Vehicle Sales
$ws.Range("C116").Value = "97274"
$ws.Range("D116").Value = "85540"
$ws.Range("H116").Value = "Oct/24"

# Business - Mining Production
$ws.Range("C119").Value = "-5.3"
$ws.Range("D119").Value = "-5.2"
$ws.Range("H119").Value = "Sep/24"

# Consumer - Retail Sales MoM
$ws.Range("C122").Value = "2.3"
$ws.Range("D122").Value = "2.3"
$ws.Range("H122").Value = "Sep/24"

# Consumer - Retail Sales YoY
$ws.Range("C123").Value = "15.9"
$ws.Range("D123").Value = "13.7"
$ws.Range("H123").Value = "Sep/24"

# Housing - Home Sales
$ws.Range("C131").Value = "165138"
$ws.Range("D131").Value = "140919"
$ws.Range("H131").Value = "Oct/24"

# Housing - Existing Home Sales
$ws.Range("C133").Value = "107459"
$ws.Range("D133").Value = "96061"
$ws.Range("H133").Value = "Oct/24"

# Housing - New Home Sales
$ws.Range("C136").Value = "57679"
$ws.Range("D136").Value = "44858"
$ws.Range("H136").Value = "Oct/24"
